# Generate Report for Handback
# Update the timestamp values recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" for the third file.
$overview.Range("G4").Value = "2016-08-30 06:48:17"

# zh-cn sheet - Handoff / Handback datetimes for the third file.
$zhcn.Range("H4").Value = "2016-08-30 06:48:13"
$zhcn.Range("K4").Value = "2016-08-30 06:48:45"

# de-de sheet - Handoff / Handback datetimes for the third file.
$dede.Range("H4").Value = "2016-08-30 06:48:17"
$dede.Range("K4").Value = "2016-08-30 06:48:51"
